# Re-apply the author's saved shape positions on slide 4 ("Beers of the
# USA" deck). Two shapes were nudged slightly (a drag in the UI):
#   - "Flowchart: Connector 5" (id 6)   -> new Left/Top (size unchanged)
#   - "Connector: Curved 22"   (id 23)  -> new Left/Top/Width/Height
#     (the curved connector glued between shapes 4 and 6; its routed
#     geometry has to be restated explicitly since nothing here
#     recomputes connector routing automatically the way PowerPoint's
#     UI would on a live drag)
#
# Values below are the point-unit (1/72") literals that round-trip
# through the COM Single-precision Left/Top/Width/Height properties to
# the exact target EMU offsets recorded in the saved OOXML.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

$flowConn5 = $s.Shapes.Item("Flowchart: Connector 5")
$flowConn5.Left = 429.3083
$flowConn5.Top  = 323.675433

$curved22 = $s.Shapes.Item("Connector: Curved 22")
$curved22.Left   = 459.06246
$curved22.Top    = 265.75222
$curved22.Width  = 68.2104
$curved22.Height = 47.6361
